# WorkflowAnalysisValidatorCounter and workflowAnalysisConditionsCounter add
#
# Adds two new worksheets - "ConditionsCounter" and "ValidatorsCounter" -
# right after the existing "PostFunctionCounter" sheet, populates them with
# the xpath/expectedNumber summary rows, and updates the selection/active
# sheet state to match (ValidatorsCounter becomes the active tab).

$wb = $excel.ActiveWorkbook

$pfc = $wb.Worksheets.Item("PostFunctionCounter")

# New sheet: ConditionsCounter, inserted immediately after PostFunctionCounter
$cc = $wb.Worksheets.Add($null, $pfc)
$cc.Name = "ConditionsCounter"

# New sheet: ValidatorsCounter, inserted immediately after ConditionsCounter
$vc = $wb.Worksheets.Add($null, $cc)
$vc.Name = "ValidatorsCounter"

# Headers
$cc.Range("A1").Value = "xpath"
$cc.Range("B1").Value = "expectedNumber"

$vc.Range("A1").Value = "xpath"
$vc.Range("B1").Value = "expectedNumber"

# Populate data rows in the same order the strings were originally authored
# (keeps shared-string table ordering/content aligned with the source edit).
$cc.Range("A3").Value = "//span[contains(., 'Eating')]/ancestor::tr/td[5]/aui-badge"
$vc.Range("A2").Value = "//span[contains(., 'Create')]/ancestor::tr/td[6]/aui-badge"
$cc.Range("A2").Value = "//span[contains(., 'Create')]/ancestor::tr/td[5]"
$vc.Range("A3").Value = "//span[contains(., 'Eating')]/ancestor::tr/td[6]"

$cc.Range("B3").Value = 1
$vc.Range("B2").Value = 1

# Column widths (matching the source: A fits the long xpath text, B fits the
# short "expectedNumber" header)
$cc.Columns("A").ColumnWidth = 105.7109375
$cc.Columns("B").ColumnWidth = 16.7109375

$vc.Columns("A").ColumnWidth = 52.7109375
$vc.Columns("B").ColumnWidth = 16.7109375

# PostFunctionCounter loses its "active tab" / single-cell selection in favor
# of a full A1:B5 selection.
$pfc.Activate()
$pfc.Range("A1:B5").Select() | Out-Null

# ConditionsCounter keeps a stray selection at A32 (left over from editing)
# but is not the active tab.
$cc.Activate()
$cc.Range("A32").Select() | Out-Null

# ValidatorsCounter ends up the active tab, selection parked on A3.
$vc.Activate()
$vc.Range("A3").Select() | Out-Null
